$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("N2").Value = 9

# Row 6
$ws.Range("G6").Value = 1.4
$ws.Range("H6").Value = 4.15
$ws.Range("I6").Value = 6.3
$ws.Range("J6").Value = 1.88
$ws.Range("K6").Value = 2.35
$ws.Range("L6").Value = 5.8
$ws.Range("P6").Value = 3.98
$ws.Range("R6").Value = 2.07
$ws.Range("W6").Value = 6.5
$ws.Range("X6").Value = 6.1
$ws.Range("Y6").Value = 7
$ws.Range("Z6").Value = 8.25
$ws.Range("AC6").Value = 13
$ws.Range("AD6").Value = 7.4
$ws.Range("AF6").Value = 50
$ws.Range("AG6").Value = 300
$ws.Range("AH6").Value = 15.5
$ws.Range("AI6").Value = 35
$ws.Range("AJ6").Value = 16
$ws.Range("AK6").Value = 100
$ws.Range("AL6").Value = 50
$ws.Range("AN6").Value = 3.3
$ws.Range("AO6").Value = 6.4
$ws.Range("AP6").Value = 15
$ws.Range("AQ6").Value = 18
$ws.Range("AT6").Value = 3.1
$ws.Range("AW6").Value = 7.8
$ws.Range("BA6").Value = 200
$ws.Range("BB6").Value = 400

# Row 8
$ws.Range("N8").Value = 13
